# Apply the commit: replace the HMP illumina/pyro test rows with the new
# OV-2 metagenome search terms, mark the "join" column as true, and drop
# the now-removed pyro rows (sheet shrinks from A1:D7 to A1:D3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ruler is hidden by default in this workbook (showRuler="0"); the saved
# copy no longer carries that explicit override, so restore the default.
$excel.ActiveWindow.DisplayRuler = $true

# Remove the now-unused rows (5,6,7 then 4) -- delete from the bottom up
# so row numbers of the rows still to delete don't shift.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Row 1 stays mostly the same, just the "join" header string is rewritten
# in the shared-strings table (same text, new string id) -- no value
# change needed here, but set it explicitly for clarity/robustness.
$ws.Range("A1").Value = "mg_f"
$ws.Range("B1").Value = "mg_r"
$ws.Range("C1").Value = "wgs_technology"
$ws.Range("D1").Value = "join"

# Row 2: first illumina search path + join = TRUE
$ws.Range("A2").Value = "/mnt/stepanauskas_nfs/ebecraft/OV2_metagenome/OV-2_P2_metagenome.fastq"
$ws.Range("B2").Value = "None"
$ws.Range("C2").Value = "illumina"
$ws.Range("D2").Value = $true

# Row 3: second illumina search path + join = TRUE
$ws.Range("A3").Value = "/mnt/stepanauskas_nfs/ebecraft/OV2_metagenome/OV-2_P3_metagenome.fastq"
$ws.Range("B3").Value = "None"
$ws.Range("C3").Value = "illumina"
$ws.Range("D3").Value = $true
